# Auto-save via app Streamlit
# 1) Normalise the phone-number cell C2 from a text string ("393316017468.0")
#    to a proper numeric value (393316017468).
# 2) Remove the cancelled booking in row 45 ("gauthier charroin"); Excel
#    shifts every row below it up by one, so the former row 46 becomes the
#    new row 45, and so on down through the former TOTAL row (51 -> 50).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 393316017468

$ws.Rows("45:45").Delete()

"done"
